# Apply "Updated questions and topics." edit to the Questions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions")

# 1. Trim the governor-limits answer: drop the leading sentence and
#    normalize the blank-line separators to single newlines.
$ws.Range("C3").Value = "They help to ensure that no one monopolizes the shared resources (Storage, CPU, Memory).`nHere are some few examples of governor limits in Salesforce:`nMaximum CPU time on the Salesforce servers - 10,000ms `nTotal number of sendEmail methods allowed - 10`nTotal number of records retrieved by a SOQL query - 50,000 `nTotal number of callouts in a transaction - 100"

# 2. Fill in the answer for "Difference between Custom Settings and Custom
#    Metadata types" (row 14) which was previously blank.
$ws.Range("C14").Value = "•Custom metadata does not support hierarchy type of data based on user profile or a specific user.`n•Custom settings data cannot be deployed using packages or metadata API/Change Sets.`n•Custom settings do not support relationship fields.`n•Custom setting data is not visible in test classes whereas metadata types are visible in test class without the “SeeAllData” annotation.`n•Custom metadata records are deployable and packageable, but Custom setting data is not."

# 3. Insert a new row right after it with a link to the source article.
$ws.Rows(15).Insert()
$ws.Range("C15").Value = "https://sfdcgenius.com/difference-between-custom-settings-and-custom-metadata-types/"
$ws.Hyperlinks.Add($ws.Range("C15"), "https://sfdcgenius.com/difference-between-custom-settings-and-custom-metadata-types/")
$ws.Rows(15).RowHeight = 30

# Keep the active selection close to where the edits were made.
[void]$ws.Range("B35").Select()
